# GIT.pptx -> "Neue Version mit Meinung"
#
# 1) Bump the cached datetimeFigureOut field text on every slide layout
#    and on the slide master from "19.01.19" to "20.01.19".
# 2) Slide 2 ("Meinungen"): merge the three runs of the second body
#    paragraph ("Now that they " / "use GIT " / "they spend so much
#    more time in source control than ever before.") into a single run.

$p = $ppt.ActivePresentation

# --- 1) Date placeholder on the slide master -----------------------------
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.Name -like "Datumsplatzhalter*") {
        if ($shp.TextFrame.TextRange.Text -eq "19.01.19") {
            $shp.TextFrame.TextRange.Text = "20.01.19"
        }
    }
}

# --- 1) Date placeholder on every slide layout ----------------------------
$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $lay = $layouts.Item($i)
    for ($j = 1; $j -le $lay.Shapes.Count; $j++) {
        $shp = $lay.Shapes.Item($j)
        if ($shp.Name -like "Datumsplatzhalter*") {
            if ($shp.TextFrame.TextRange.Text -eq "19.01.19") {
                $shp.TextFrame.TextRange.Text = "20.01.19"
            }
        }
    }
}

# --- 2) Merge the runs of the second paragraph on slide 2 ----------------
$s2 = $p.Slides.Item(2)
$shape = $s2.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

$target = "Now that they use GIT they spend so much more time in source control than ever before."

$paraCount = $tr.Paragraphs().Count
for ($k = 1; $k -le $paraCount; $k++) {
    $para = $tr.Paragraphs($k, 1)
    $paraText = $para.Text.TrimEnd("`r")
    if ($paraText -eq $target) {
        $prev = $tr.Paragraphs($k - 1, 1)
        $para.Delete()
        $inserted = $prev.InsertAfter("`r" + $target)
        break
    }
}
